# Implement search by cancer type.
#
# Adds a new "Basic Cancer Types" worksheet between "BasicSearch" and
# "AdvanceSearch", populates it with lookup data, and updates the active
# selection on the workbook / BasicSearch sheet accordingly.

$wb = $excel.ActiveWorkbook

$basic = $wb.Worksheets.Item("BasicSearch")

# BasicSearch's own selection moves from E4 to E1 once it stops being the
# active tab.
$null = $basic.Range("E1").Select()

# Insert the new sheet right after BasicSearch (i.e. before AdvanceSearch).
$cancerTypes = $wb.Worksheets.Add($null, $basic)
$cancerTypes.Name = "Basic Cancer Types"

# Header row - copy the header formatting used on BasicSearch's row 1.
$cancerTypes.Range("A1").Value = "Cancer Type Name"
$cancerTypes.Range("B1").Value = "Concept ID"
$null = $basic.Range("A1").Copy()
$null = $cancerTypes.Range("A1:B1").PasteSpecial(-4122)

# Data rows.
$cancerTypes.Range("A2").Value = "Adenosquamous Lung Cancer"
$cancerTypes.Range("B2").Value = "C9133"
$cancerTypes.Range("A3").Value = "Chronic Kidney Disease, Stage 2"
$cancerTypes.Range("B3").Value = "C80388"
$cancerTypes.Range("A4").Value = "Childhood Giant Cell Glioblastoma"
$cancerTypes.Range("B4").Value = "C114966"
$cancerTypes.Range("A5").Value = "Chronic Atrophic Gastritis"
$cancerTypes.Range("B5").Value = "C7405"

# Column widths.
$cancerTypes.Columns.Item(1).ColumnWidth = 35.7109375
$cancerTypes.Columns.Item(2).ColumnWidth = 12.85546875

# Selection / active sheet.
$null = $cancerTypes.Range("A6").Select()
$null = $cancerTypes.Select()
